$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the columns C:F entirely (no longer part of the table)
$ws.Range("C1:F6").EntireColumn.Delete()

# Remove rows 5:6 entirely (table is now only 4 rows tall)
$ws.Range("A5:A6").EntireRow.Delete()

# Update header row (A1 should match B1's existing bold/border header style)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Student_id"
$ws.Range("B1").Value = "Student_name"

# Update data rows - these plain data cells should carry no special style
$ws.Range("A2").Value = 205
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Ирина"

$ws.Range("A3").Value = 206
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Ксения"

$ws.Range("A4").Value = 207
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "Жанна"
